$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 45175
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = 100112001
$ws.Range("G96").Value = "Berenjena"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 150
$ws.Range("K96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = 10000
$ws.Range("N96").Value = "$/caja 50 unidades"
$ws.Range("O96").Value = "Región de Arica y Parinacota"
$ws.Range("P96").Value = 200
$ws.Range("Q96").Value = 50
$ws.Range("R96").Value = "Hortaliza"
